# Changes to getUserGroupNotifications: added additional filter groupName
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("API's V2")

$purpose = "To get all user notifications `n`"notificationType`", `"groupName`" are  optional.`nExpected values are TEXT, FILE`nFILE : To get only File Notifications`nTEXT : To get only Text Notifications"
$request = "{`n    `"email`": `"notificationboard1@gmail.com`",`n    `"notificationType`": `"FILE`",`n    `"groupName`": `"New`"`n}"

$ws.Range("D8").Value = $request
$ws.Range("C8").Value = $purpose

$ws.Activate()
$ws.Range("C8").Select()
